$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Entry Condition cell: append a new sentence after the existing
#    "... deve trovarsi nella pagina del carrello" text.
# ------------------------------------------------------------------
$entry = $d.Content
$entry.Find.Execute("deve trovarsi nella pagina del carrello", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$entry.Collapse(0)
$entry.InsertAfter(" e deve essere presente almeno un prodotto nel carrello")

# ------------------------------------------------------------------
# 2) Re-create the "_GoBack" bookmark at the end of that same
#    paragraph (this is where Word leaves it after the last edit).
#    Adding it here bumps every other bookmark's id by one and
#    removes the previous "_GoBack" bookmark automatically.
# ------------------------------------------------------------------
$entry2 = $d.Content
$entry2.Find.Execute("e deve essere presente almeno un prodotto nel carrello", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$entry2.Collapse(0)
Write-Output ("DEBUG bookmark target range: [" + $entry2.Start + "," + $entry2.End + "]")
$d.Bookmarks.Add("_GoBack", $entry2)

# ------------------------------------------------------------------
# 3) Merge the two runs that make up the failure-scenario heading
#    "Scenario/Flusso di eventi in caso di fallimento: La quantità
#    dei prodotti richiesti dall'utente eccede la disponibilità"
#    into a single run.
# ------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("Scenario/Flusso di eventi in caso di fallimento: La quantità dei prodotti richiesti dall’utente eccede la disponibilità", `
    $false, $true, $false, $false, $false, $true, 1, $false, `
    "Scenario/Flusso di eventi in caso di fallimento: La quantità dei prodotti richiesti dall’utente eccede la disponibilità", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Merge the two runs "Visualizza un messaggio che riferisce
#    all'utente che " + "non vi è disponibilità sufficiente" into a
#    single run.
# ------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Visualizza un messaggio che riferisce all’utente che non vi è disponibilità sufficiente", `
    $false, $true, $false, $false, $false, $true, 1, $false, `
    "Visualizza un messaggio che riferisce all’utente che non vi è disponibilità sufficiente", 2) | Out-Null
